# Trans-Tasman charts updated to 29 November
# Insert 7 new days of VIC second-dose data at the top (rows 2-8), shifting
# the rest of the table down, update the selection, header font colour and
# column widths to match the refreshed export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 7 new rows above the existing data (row 2), pushing everything
#    else down by 7 rows.
$ws.Rows.Item(2).Resize(7).Insert()

# 1b. Clone the formatting (number format + font) of the first data row
#     (now row 9) onto the freshly inserted rows, using a format-only
#     paste so the shared style table isn't forked.
$ws.Range("A9:B9").Copy()
$ws.Range("A2:B8").PasteSpecial(-4122)
$ws.Rows.Item(2).Resize(7).RowHeight = 18

# 2. New data rows (most-recent-first).
$newDates = @(44529, 44528, 44527, 44526, 44525, 44524, 44523)
$newDoses = @(5140389, 5138001, 5132289, 5118977, 5105511, 5092441, 5080634)

for ($i = 0; $i -lt 7; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $newDates[$i]
    $ws.Cells.Item($r, 2).Value = $newDoses[$i]
}

# 3. Header font turns red (bold red Arial 9.6 on B1).
$ws.Range("B1").Font.Color = 255

# 4. Column widths: column A keeps its best-fit width, column B gets its
#    own (slightly narrower) best-fit width.
$ws.Columns.Item(1).ColumnWidth = 12.6640625
$ws.Columns.Item(2).ColumnWidth = 11.83203125

# 5. Selection moves to B1.
$ws.Range("B1").Select()
